$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

$replacements = @(
    @("2024-10-14 Monday", "2024-10-15 Tuesday"),
    @("23×16=", "31×37="),
    @("14×64=", "76×26="),
    @("66×89=", "94×52="),
    @("99×74=", "20×61="),
    @("84×65=", "73×48="),
    @("45×98=", "29×40="),
    @("55×65=", "39×23="),
    @("95×38=", "47×79="),
    @("13×98=", "47×86="),
    @("24×68=", "87×25="),
    @("14×50=", "89×31="),
    @("93×87=", "14×58="),
    @("63×73=", "82×40="),
    @("62×44=", "99×29="),
    @("62×89=", "77×83="),
    @("80×73=", "55×23="),
    @("34×37=", "76×73="),
    @("82×44=", "74×37="),
    @("95×89=", "52×39="),
    @("25×54=", "88×95="),
    @("68×93=", "45×68="),
    @("94×23=", "31×59="),
    @("98×50=", "76×28="),
    @("80×44=", "74×11="),
    @("92×74=", "80×49=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, $wdReplaceAll)
}
